$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "Figure 3.2"
$ws.Range("B9").Value = "Results and Discussion"
$ws.Range("C9").Value = "Evolutionary history of phototransduction components gene families and distribution across Eukarya"
$ws.Range("D9").Value = "Added to manuscript."

$ws.Range("C16").Select()
